$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "quality_comparison"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

# B4 already carries the plain "no border / default font / default alignment"
# style -- use it as a clean format donor so C1/D1 drop the bold+centered
# header look and fall back to the worksheet default font before the new
# borders go on.
$ws1.Range("B4").Copy()
$ws1.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("D1").PasteSpecial(-4122)   # xlPasteFormats

# C1: top+bottom-only thin border (matches border id 4).
$c1 = $ws1.Range("C1")
$c1.Borders.LineStyle = 1        # xlContinuous on every edge...
$c1.Borders.Item(7).LineStyle = -4142   # ...then drop xlEdgeLeft
$c1.Borders.Item(10).LineStyle = -4142  # ...and xlEdgeRight

# D1: top+bottom+right thin border (matches border id 5).
$d1 = $ws1.Range("D1")
$d1.Borders.LineStyle = 1        # xlContinuous on every edge...
$d1.Borders.Item(7).LineStyle = -4142   # ...then drop xlEdgeLeft only

# Anonymize the "fedcore" header label.
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------------
# Sheet "computational_comparison"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

$ws2.Range("B4").Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("D1").PasteSpecial(-4122)

$c1b = $ws2.Range("C1")
$c1b.Borders.LineStyle = 1
$c1b.Borders.Item(7).LineStyle = -4142
$c1b.Borders.Item(10).LineStyle = -4142

$d1b = $ws2.Range("D1")
$d1b.Borders.LineStyle = 1
$d1b.Borders.Item(7).LineStyle = -4142

# F1/G1 need exactly the same two resulting styles as C1/D1 -- clone the
# already-finished formats instead of re-deriving them from scratch so the
# style table doesn't pick up a duplicate/transient entry.
$c1b.Copy()
$ws2.Range("F1").PasteSpecial(-4122)

$d1b.Copy()
$ws2.Range("G1").PasteSpecial(-4122)

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5.
$ws2.Range("G5").ClearContents()
